# Generate Report for Handback
# Renames the two handed-back files (new GUIDs) and refreshes their
# handoff/handback timestamps & xliff file names across all three sheets.

$wb = $excel.ActiveWorkbook

$oldGuid1 = "1b31e8cc-026b-48c4-b541-d3126b3474cc"
$newGuid1 = "2f6ddd41-cb43-4754-ab4c-4f00772dd94c"
$oldGuid2 = "8c8f818e-736a-455b-af0a-3d9d4145e437"
$newGuid2 = "ffff3cbbd39e-7737-42d0-8a88-675e3972ed4a"

$newHash = "8c4add3d4810929d8c9d72d15b02db7b90fc5767"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid1.md"
$wsOverview.Range("C2").Value = ".md"
$wsOverview.Range("G2").Value = "2016-09-03 21:05:45"

$wsOverview.Range("A3").Value = "$newGuid2.md"
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("G3").Value = "2016-09-03 21:05:45"

$wsOverview.Range("B2").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/36112e37e0d051197d8101660f263f62debb5430/e2e/$oldGuid1.md", "", "", "e2e\$newGuid1.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/36112e37e0d051197d8101660f263f62debb5430/e2e/$oldGuid2.md", "", "", "e2e\$newGuid2.md")

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$xlfZhCn = "$newGuid1.$newHash.zh-cn.xlf"

$wsZhCn.Range("A2").Value = "$newGuid1.md"
$wsZhCn.Range("G2").Value = $xlfZhCn
$wsZhCn.Range("H2").Value = "2016-09-03 21:05:39"
$wsZhCn.Range("I2").Value = "$newGuid1.md"
$wsZhCn.Range("J2").Value = $xlfZhCn
$wsZhCn.Range("K2").Value = "2016-09-03 21:05:56"

$wsZhCn.Range("A3").Value = "$newGuid2.md"
$wsZhCn.Range("G3").Value = $xlfZhCn
$wsZhCn.Range("H3").Value = "2016-09-03 21:05:39"
$wsZhCn.Range("I3").Value = "$newGuid2.md"
$wsZhCn.Range("J3").Value = $xlfZhCn
$wsZhCn.Range("K3").Value = "2016-09-03 21:05:56"

$wsZhCn.Range("A2").Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/36112e37e0d051197d8101660f263f62debb5430/e2e/$oldGuid1.md", "", "", "$newGuid1.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/e7ba5ad459badf4093af1c1483d860cd7629cae6/e2e/$oldGuid1.md", "", "", "$newGuid1.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/36112e37e0d051197d8101660f263f62debb5430/e2e/$oldGuid2.md", "", "", "$newGuid2.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/e7ba5ad459badf4093af1c1483d860cd7629cae6/e2e/$oldGuid2.md", "", "", "$newGuid2.md")

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$xlfDeDe = "$newGuid1.$newHash.de-de.xlf"

$wsDeDe.Range("A2").Value = "$newGuid1.md"
$wsDeDe.Range("G2").Value = $xlfDeDe
$wsDeDe.Range("H2").Value = "2016-09-03 21:05:45"
$wsDeDe.Range("I2").Value = "$newGuid1.md"
$wsDeDe.Range("J2").Value = $xlfDeDe
$wsDeDe.Range("K2").Value = "2016-09-03 21:06:09"

$wsDeDe.Range("A3").Value = "$newGuid2.md"
$wsDeDe.Range("G3").Value = $xlfDeDe
$wsDeDe.Range("H3").Value = "2016-09-03 21:05:45"
$wsDeDe.Range("I3").Value = "$newGuid2.md"
$wsDeDe.Range("J3").Value = $xlfDeDe
$wsDeDe.Range("K3").Value = "2016-09-03 21:06:09"

$wsDeDe.Range("A2").Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/36112e37e0d051197d8101660f263f62debb5430/e2e/$oldGuid1.md", "", "", "$newGuid1.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/dd1a6d9d71d7c1cde94387d4d5ac2c254726687a/e2e/$oldGuid1.md", "", "", "$newGuid1.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/36112e37e0d051197d8101660f263f62debb5430/e2e/$oldGuid2.md", "", "", "$newGuid2.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/dd1a6d9d71d7c1cde94387d4d5ac2c254726687a/e2e/$oldGuid2.md", "", "", "$newGuid2.md")

Write-Host "Handback status report regenerated."
